$d = $word.ActiveDocument

function Insert-ParagraphXml($range, $innerXml) {
    $pkg = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# --- Edit 1: "Amount of completed tasks:" paragraph gets " 10" appended
# (bold run with a space, then a proofErr-wrapped bold run "10").
$pTasks = $d.Paragraphs(3)
$rTasks = $pTasks.Range
$selTasks = $d.Range($rTasks.Start, $rTasks.End - 1)
$xmlTasks = '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="00970715"><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Amount of completed tasks:</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>10</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '</w:p>'
Insert-ParagraphXml $selTasks $xmlTasks

# --- Edit 2: "...ndone or incomplete:" paragraph gets " 0" appended
# (bold run with a space, then a proofErr-wrapped bold run "0").
$pUndone = $d.Paragraphs(4)
$rUndone = $pUndone.Range
$selUndone = $d.Range($rUndone.Start, $rUndone.End - 1)
$xmlUndone = '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="00970715"><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Which tasks were left u</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>ndone or incomplete:</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>0</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '</w:p>'
Insert-ParagraphXml $selUndone $xmlUndone

# --- Edit 3: the empty "Otsikko2"-styled paragraph (right after
# "Self-assessment: ") becomes a plain paragraph with red self-assessment
# text, followed by a new empty plain paragraph (no pStyle) before the
# "Answers to other than coding tasks here:" paragraph.
$pEmpty = $d.Paragraphs(6)
$rEmpty = $pEmpty.Range
$selEmpty = $d.Range($rEmpty.Start, $rEmpty.End)
$xmlEmpty = '<w:p><w:pPr><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-US"/></w:rPr><w:t>I did a good job with most of the exercises except for the last one, which was a struggle.</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
Insert-ParagraphXml $selEmpty $xmlEmpty

Write-Host "edits applied"
